$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell F1 (copy formatting from E1, the last header cell, then set its text)
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Add time_taken values for each data row as plain text
$ws.Range("F2").Value = "2021-10-05 13:39:44.821495"
$ws.Range("F3").Value = "2021-10-05 13:39:44.821505"
$ws.Range("F4").Value = "2021-10-05 13:39:44.821509"
$ws.Range("F5").Value = "2021-10-05 13:39:44.821512"
$ws.Range("F6").Value = "2021-10-05 13:39:44.821514"
$ws.Range("F7").Value = "2021-10-05 13:39:44.821517"
$ws.Range("F8").Value = "2021-10-05 13:39:44.821520"
$ws.Range("F9").Value = "2021-10-05 13:39:44.821522"
$ws.Range("F10").Value = "2021-10-05 13:39:44.821525"
$ws.Range("F11").Value = "2021-10-05 13:39:44.821528"
